# Applies the Mon Nov  6 14:23:22 UTC 2023 "cryptos list" refresh:
#  - Updates Price (D) and Volume(1h) (E) figures for the existing rows.
#  - A handful of rows reshuffled rank order, which is reflected here as
#    updates to Coin (B) / Link (C) / Price (D) / Volume(1h) (E) in place,
#    leaving the rank column (A) untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds values as plain text (many contain more than one
# "." thousands separator, e.g. "35.396.90"), so force Text format before
# writing the new figures to stop Excel from re-interpreting the plain
# numeric-looking ones (e.g. "0.738") as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '35.396.90'
$ws.Range("E2").Value = '  +0.12%  '

# Row 3
$ws.Range("D3").Value = '1.921.58'
$ws.Range("E3").Value = '  +1.19%  '

# Row 4
$ws.Range("E4").Value = '  -0.44%  '

# Row 5
$ws.Range("D5").Value = '0.738'
$ws.Range("E5").Value = '  +13.35%  '

# Row 6
$ws.Range("D6").Value = '253.13'
$ws.Range("E6").Value = '  +3.83%  '

# Row 7
$ws.Range("E7").Value = '  -0.48%  '

# Row 8
$ws.Range("D8").Value = '40.72'
$ws.Range("E8").Value = '  -1.90%  '

# Row 9
$ws.Range("D9").Value = '0.356'
$ws.Range("E9").Value = '  +4.41%  '

# Row 10
$ws.Range("D10").Value = '52.57'
$ws.Range("E10").Value = '  +5.29%  '

# Row 11
$ws.Range("D11").Value = '0.0738'
$ws.Range("E11").Value = '  +3.92%  '

# Row 12
$ws.Range("D12").Value = '0.0998'
$ws.Range("E12").Value = '  +0.15%  '

# Row 13
$ws.Range("D13").Value = '2.198.77'
$ws.Range("E13").Value = '  +1.05%  '

# Row 14
$ws.Range("D14").Value = '12.66'
$ws.Range("E14").Value = '  +4.50%  '

# Row 15
$ws.Range("D15").Value = '0.716'
$ws.Range("E15").Value = '  +3.33%  '

# Row 16
$ws.Range("D16").Value = '1.917.60'
$ws.Range("E16").Value = '  +0.89%  '

# Row 17
$ws.Range("E17").Value = '  +1.14%  '

# Row 18
$ws.Range("D18").Value = '35.394.97'
$ws.Range("E18").Value = '  +0.06%  '

# Row 19
$ws.Range("D19").Value = '73.27'
$ws.Range("E19").Value = '  +2.35%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0834'
$ws.Range("E20").Value = '  +2.37%  '

# Row 21
$ws.Range("D21").Value = '13.05'
$ws.Range("E21").Value = '  +4.31%  '

# Row 22
$ws.Range("D22").Value = '241.80'
$ws.Range("E22").Value = '  +0.03%  '

# Row 23
$ws.Range("D23").Value = '5.08'
$ws.Range("E23").Value = '  +7.45%  '

# Row 24
$ws.Range("E24").Value = '  -0.53%  '

# Row 25
$ws.Range("D25").Value = '2.37'
$ws.Range("E25").Value = '  +2.94%  '

# Row 26
$ws.Range("D26").Value = '2.33'
$ws.Range("E26").Value = '  +0.56%  '

# Row 27
$ws.Range("D27").Value = '167.77'
$ws.Range("E27").Value = '  -1.43%  '

# Row 28
$ws.Range("D28").Value = '8.72'
$ws.Range("E28").Value = '  +4.63%  '

# Row 29
$ws.Range("D29").Value = '0.136'
$ws.Range("E29").Value = '  +7.52%  '

# Row 30
$ws.Range("D30").Value = '18.83'
$ws.Range("E30").Value = '  +3.38%  '

# Row 31
$ws.Range("D31").Value = '4.131.87'
$ws.Range("E31").Value = '  +19.55%  '

# Row 32 (now TrustWalletToken)
$ws.Range("B32").Value = 'TrustWalletToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D32").Value = '1.68'
$ws.Range("E32").Value = '  +26.74%  '

# Row 33 (now Filecoin)
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '4.35'
$ws.Range("E33").Value = '  +5.50%  '

# Row 34 (now WEMIXToken)
$ws.Range("B34").Value = 'WEMIXToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").Value = '1.99'
$ws.Range("E34").Value = '  +14.32%  '

# Row 35
$ws.Range("D35").Value = '0.0581'
$ws.Range("E35").Value = '  +3.36%  '

# Row 36
$ws.Range("D36").Value = '4.26'
$ws.Range("E36").Value = '  +4.15%  '

# Row 37
$ws.Range("E37").Value = '  -0.42%  '

# Row 38
$ws.Range("E38").Value = '  -2.40%  '

# Row 39
$ws.Range("D39").Value = '2.03'
$ws.Range("E39").Value = '  -0.39%  '

# Row 40
$ws.Range("D40").Value = '17.47'
$ws.Range("E40").Value = '  +10.32%  '

# Row 41
$ws.Range("D41").Value = '99.05'
$ws.Range("E41").Value = '  +10.90%  '

# Row 42
$ws.Range("E42").Value = '  +4.90%  '

# Row 43
$ws.Range("E43").Value = '  -0.37%  '

# Row 44
$ws.Range("D44").Value = '0.0653'
$ws.Range("E44").Value = '  +3.72%  '

# Row 45 (now RenderToken)
$ws.Range("B45").Value = 'RenderToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D45").Value = '2.49'
$ws.Range("E45").Value = '  +5.94%  '

# Row 46 (now Maker)
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.347.76'
$ws.Range("E46").Value = '  +0.79%  '

# Row 47
$ws.Range("E47").Value = '  +0.51%  '

# Row 48 (now MXToken)
$ws.Range("B48").Value = 'MXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D48").Value = '2.78'
$ws.Range("E48").Value = '  +0.10%  '

# Row 49 (now FraxShare)
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '6.71'
$ws.Range("E49").Value = '  +3.22%  '

# Row 50
$ws.Range("D50").Value = '45.32'
$ws.Range("E50").Value = '  -4.92%  '

# Row 51 (now RocketPoolETH)
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.105.98'
$ws.Range("E51").Value = '  +1.01%  '
